$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1466.25
$ws.Range("J112").Value = 1511.0938
$ws.Range("L112").Value = 4533.2814
$ws.Range("N112").Value = -6749.2814

# Row 138
$ws.Range("H138").Value = 3193.2654
$ws.Range("I138").Value = 2081.5833
$ws.Range("J138").Value = 3553.8108
$ws.Range("K138").Value = 6244.749899999999
$ws.Range("L138").Value = 10661.4324
$ws.Range("M138").Value = -1104.749899999999
$ws.Range("N138").Value = -20941.4324

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 4410.2104
$ws.Range("I63").Value = 2566.1667
$ws.Range("K63").Value = 2566.1667
$ws.Range("M63").Value = -1880.1667

# Row 66
$ws.Range("H66").Value = 4410.2104
$ws.Range("I66").Value = 2566.1667
$ws.Range("K66").Value = 12830.8335
$ws.Range("M66").Value = -9398.833500000001

# Row 132
$ws.Range("H132").Value = 2325.9167
$ws.Range("I132").Value = 1942.9512
$ws.Range("K132").Value = 5828.8536
$ws.Range("M132").Value = -3298.8536

# Row 139
$ws.Range("H139").Value = 79759.375
$ws.Range("J139").Value = 79759.375
$ws.Range("L139").Value = 79759.375
$ws.Range("N139").Value = -90039.375

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 701.4167
$ws.Range("I64").Value = 684.6
$ws.Range("J64").Value = 713.4286
$ws.Range("K64").Value = 684.6
$ws.Range("L64").Value = 713.4286
$ws.Range("M64").Value = -459.6
$ws.Range("N64").Value = -1163.4286

# Row 67
$ws.Range("H67").Value = 701.4167
$ws.Range("I67").Value = 684.6
$ws.Range("J67").Value = 713.4286
$ws.Range("K67").Value = 684.6
$ws.Range("L67").Value = 713.4286
$ws.Range("M67").Value = 95.39999999999998
$ws.Range("N67").Value = -2273.4286

# Row 94
$ws.Range("H94").Value = 31251340
$ws.Range("I94").Value = 38462496
$ws.Range("K94").Value = 38462496
$ws.Range("M94").Value = -38462045

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1500
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2700

# Row 132
$ws.Range("H132").Value = 3512.1667
$ws.Range("I132").Value = 3380
$ws.Range("K132").Value = 10140
$ws.Range("M132").Value = -7610

# Row 134
$ws.Range("H134").Value = 2389.2546
$ws.Range("I134").Value = 2012.66
$ws.Range("K134").Value = 6037.98
$ws.Range("M134").Value = -3502.98

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 2316.7
$ws.Range("I14").Value = 2316.7
$ws.Range("K14").Value = 6950.099999999999
$ws.Range("M14").Value = -6777.099999999999

# Row 86
$ws.Range("H86").Value = 600.5
$ws.Range("J86").Value = 767.3333
$ws.Range("L86").Value = 2301.9999
$ws.Range("N86").Value = -4673.9999

# Row 89
$ws.Range("H89").Value = 600.5
$ws.Range("J89").Value = 767.3333
$ws.Range("L89").Value = 6905.9997
$ws.Range("N89").Value = -18761.9997

# Row 113
$ws.Range("H113").Value = 1821.3572
$ws.Range("J113").Value = 1880.5769
$ws.Range("L113").Value = 5641.7307
$ws.Range("N113").Value = -9981.7307

# Row 137
$ws.Range("H137").Value = 4350881
$ws.Range("I137").Value = 11112634
$ws.Range("J137").Value = 4040
$ws.Range("K137").Value = 33337902
$ws.Range("L137").Value = 12120
$ws.Range("M137").Value = -33332802
$ws.Range("N137").Value = -22320

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 958.5517
$ws.Range("I2").Value = 592.3333
$ws.Range("J2").Value = 1350.9286
$ws.Range("K2").Value = 592.3333
$ws.Range("L2").Value = 1350.9286
$ws.Range("M2").Value = -479.3333
$ws.Range("N2").Value = -1576.9286

# Row 43
$ws.Range("H43").Value = 20803
$ws.Range("I43").Value = 1338.6666
$ws.Range("J43").Value = 49999.5
$ws.Range("K43").Value = 1338.6666
$ws.Range("L43").Value = 49999.5
$ws.Range("M43").Value = -1187.6666
$ws.Range("N43").Value = -50301.5

# Row 46
$ws.Range("H46").Value = 38553.31
$ws.Range("J46").Value = 44998.25
$ws.Range("L46").Value = 44998.25
$ws.Range("N46").Value = -45310.25

# Row 57
$ws.Range("H57").Value = 17568.625
$ws.Range("J57").Value = 26666.334
$ws.Range("L57").Value = 26666.334
$ws.Range("N57").Value = -28306.334

# Row 80
$ws.Range("H80").Value = 3301.1765
$ws.Range("I80").Value = 2661.5833
$ws.Range("J80").Value = 4836.2
$ws.Range("K80").Value = 2661.5833
$ws.Range("L80").Value = 4836.2
$ws.Range("M80").Value = -1663.5833
$ws.Range("N80").Value = -6832.2

# Row 83
$ws.Range("H83").Value = 3301.1765
$ws.Range("I83").Value = 2661.5833
$ws.Range("J83").Value = 4836.2
$ws.Range("K83").Value = 13307.9165
$ws.Range("L83").Value = 24181
$ws.Range("M83").Value = -8315.916499999999
$ws.Range("N83").Value = -34165

# Row 102
$ws.Range("H102").Value = 1342.9231
$ws.Range("I102").Value = 1422.5834
$ws.Range("K102").Value = 1422.5834
$ws.Range("M102").Value = 199.4166

# Row 107
$ws.Range("H107").Value = 35715184
$ws.Range("I107").Value = 185
$ws.Range("J107").Value = 100002180
$ws.Range("K107").Value = 185
$ws.Range("L107").Value = 100002180
$ws.Range("M107").Value = 1735
$ws.Range("N107").Value = -100006020

# Row 132
$ws.Range("H132").Value = 4821.9033
$ws.Range("I132").Value = 3611.3809
$ws.Range("K132").Value = 10834.1427
$ws.Range("M132").Value = -8304.1427

# Row 134
$ws.Range("H134").Value = 109988
$ws.Range("J134").Value = 109988
$ws.Range("L134").Value = 329964
$ws.Range("N134").Value = -335034

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 13078.484
$ws.Range("I40").Value = 6316.4585
$ws.Range("K40").Value = 6316.4585
$ws.Range("M40").Value = -6180.4585

# Row 82
$ws.Range("H82").Value = 4447.207
$ws.Range("I82").Value = 3121.5833
$ws.Range("K82").Value = 3121.5833
$ws.Range("M82").Value = -2760.5833

# Row 85
$ws.Range("H85").Value = 4447.207
$ws.Range("I85").Value = 3121.5833
$ws.Range("K85").Value = 3121.5833
$ws.Range("M85").Value = -1873.5833

# Row 122
$ws.Range("H122").Value = 75641.14
$ws.Range("I122").Value = 103884.85
$ws.Range("J122").Value = 5031.875
$ws.Range("K122").Value = 311654.55
$ws.Range("L122").Value = 15095.625
$ws.Range("M122").Value = -309204.55
$ws.Range("N122").Value = -19995.625

# Row 136
$ws.Range("H136").Value = 5769.8223
$ws.Range("I136").Value = 5478.0625
$ws.Range("J136").Value = 6488
$ws.Range("K136").Value = 16434.1875
$ws.Range("L136").Value = 19464
$ws.Range("M136").Value = -13884.1875
$ws.Range("N136").Value = -24564

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 44677.145
$ws.Range("I2").Value = 44677.145
$ws.Range("K2").Value = 44677.145
$ws.Range("M2").Value = -44565.145

# Row 107
$ws.Range("H107").Value = 33335052
$ws.Range("I107").Value = 1860.4445
$ws.Range("K107").Value = 5581.333500000001
$ws.Range("M107").Value = -3661.333500000001
